$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update values B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: remove D2, add B2 and C2
$ws.Range("D2").ClearContents()
$ws.Range("B2").Value = 41.226569807504134
$ws.Range("C2").Value = 21.066300986616007

# Row 3: remove B3, update C3
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 40.31666820574489

# Update selection to match new active range
$null = $ws.Range("B1:E3").Select()
